$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.08%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.67%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.68%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07366"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.47%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.293"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "27.66%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.664"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.61%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9167"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.64%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "14.59%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1703"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.60%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08310"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.66%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03115"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.31%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.24%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001499"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.55%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005750"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.471"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.23%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.738"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.17%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.61%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3329"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.29%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1299"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.13%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.171"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.51%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.15%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04504"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001214"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004198"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.62%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.10%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003391"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.26%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01571"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04509"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.58%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009847"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.66%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.44%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002216"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.83%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008523"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.03%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006103"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.606"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "15.58%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001998"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-30.83%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
